# "write some new for stephen"
#
# The target revision carries a second, newly-inserted "Subtitle"
# paragraph style definition in styles.xml (the document already had one
# Subtitle style; the edit adds another block with the same name/
# formatting — basedOn Normal, next Normal, keepNext/keepLines, no page
# break before, 18pt space-before / 4pt space-after with auto line
# spacing, Georgia rFonts (ascii/eastAsia/hAnsi/cs), italic, color
# 666666, 24pt (sz/szCs 48)).
#
# Word's object model can't mint a second style sharing the exact same
# internal styleId (styles are looked up/created by a single Add call
# keyed on name), so we add a new style and drive it to look exactly
# like the existing "Subtitle" style: same base/next style, same
# paragraph formatting, same run formatting, and the same display name.

$d = $word.ActiveDocument

$new = $d.Styles.Add("Subtitle (new)", 1)

$new.BaseStyle = "Normal"
$new.NextParagraphStyle = "Normal"

$new.ParagraphFormat.KeepWithNext = $true
$new.ParagraphFormat.KeepTogether = $true
$new.ParagraphFormat.PageBreakBefore = $false
$new.ParagraphFormat.SpaceBefore = 18
$new.ParagraphFormat.SpaceAfter = 4
$new.ParagraphFormat.LineSpacingRule = 5

$new.Font.Name = "Georgia"
$new.Font.NameFarEast = "Georgia"
$new.Font.NameBi = "Georgia"
$new.Font.Italic = $true
$new.Font.Color = 6710886
$new.Font.Size = 24
$new.Font.SizeBi = 24

$new.NameLocal = "Subtitle"

Write-Output "done"
